# Apply updated cryptocurrency price/volume figures to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new text value (values must remain plain text,
# matching the original inline-string cell type used in the sheet)
$updates = [ordered]@{
    "D2" = "301.21"
    "E2" = "0.61%"
    "D3" = "32.17"
    "E3" = "1.94%"
    "D4" = "4.985"
    "E4" = "-3.12%"
    "D5" = "0.07892"
    "E5" = "-2.48%"
    "D6" = "2.113"
    "E6" = "-14.71%"
    "D7" = "7.797"
    "E7" = "-0.03%"
    "E8" = "-1.61%"
    "D9" = "0.9251"
    "E9" = "-0.49%"
    "E10" = "-0.87%"
    "D11" = "0.07979"
    "E11" = "7.68%"
    "D12" = "0.08727"
    "E12" = "-1.75%"
    "D13" = "0.03102"
    "E13" = "3.48%"
    "E14" = "0.34%"
    "D15" = "0.001510"
    "E15" = "0.11%"
    "D16" = "0.005822"
    "E16" = "-2.28%"
    "E17" = "2,101.68%"
    "D18" = "3.462"
    "E18" = "-1.91%"
    "D19" = "2.255"
    "E19" = "-1.50%"
    "E20" = "0.47%"
    "E21" = "-2.26%"
    "D22" = "4.295"
    "E22" = "3.29%"
    "E23" = "6.72%"
    "D24" = "0.04604"
    "D25" = "0.001236"
    "E25" = "-0.22%"
    "D26" = "0.004439"
    "E26" = "-2.12%"
    "E27" = "4.34%"
    "D39" = "0.01712"
    "E39" = "-2.41%"
    "D40" = "0.04773"
    "E40" = "3.93%"
    "D41" = "0.007437"
    "E41" = "7.49%"
    "D42" = "0.1355"
    "E42" = "-1.28%"
    "D43" = "0.002361"
    "E43" = "7.94%"
    "D44" = "0.01127"
    "E44" = "9.54%"
    "D45" = "0.00006017"
    "E45" = "-2.08%"
    "E46" = "0.15%"
    "D47" = "0.003393"
    "E47" = "-59.57%"
    "D48" = "0.8204"
    "E48" = "9.61%"
    "E49" = "0.15%"
    "E50" = "0.15%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text number format so numeric-looking strings (e.g. "301.21")
    # and percentages (e.g. "0.61%") are stored as text, not converted numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}

Write-Output "Updated $($updates.Count) cells"
